# Remove the "Transition from API Endpoints to DB" slide (3rd slide,
# sldId 262) from the deck. PowerPoint automatically renumbers the
# relationship ids / notesMaster id on save; that is a side effect of
# the host application and not something this script needs to manage
# by hand.
$p = $ppt.ActivePresentation
$p.Slides.Item(3).Delete()
